$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update the expected-message text for the "username only" and "password only" rows
$ws1.Range("C2").Value = "Please check your password"
$ws1.Range("C3").Value = "Please check your user id"

# Remove the old row 4 (username/password/Invalid Username and Password),
# shifting the remaining rows up
$ws1.Rows.Item(4).Select() | Out-Null
$ws1.Rows.Item(4).Delete() | Out-Null
